# TOPSIS lab-4 data entry: replace the generated decimal criteria sample
# data (columns C:E, rows 2-21) with freshly typed-in integer scores, and
# normalize all three criteria columns onto a single "whole number" display
# format (previously C used 0.00 and D:E used a custom 0.0 format).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New integer values for criteria columns C (Kryterium 1), D (Kryterium 2),
# E (Kryterium 3) for alternatives in rows 2-21.
$data = @{
    2  = @(94, 82, 73)
    3  = @(46, 64, 80)
    4  = @(53, 83, 51)
    5  = @(88, 79, 72)
    6  = @(80, 80, 7)
    7  = @(18, 71, 35)
    8  = @(95, 39, 63)
    9  = @(94, 37, 51)
    10 = @(33, 16, 88)
    11 = @(36, 74, 41)
    12 = @(18, 19, 46)
    13 = @(12, 67, 11)
    14 = @(11, 34, 98)
    15 = @(81, 64, 78)
    16 = @(75, 29, 3)
    17 = @(50, 36, 35)
    18 = @(96, 30, 96)
    19 = @(33, 82, 8)
    20 = @(64, 11, 15)
    21 = @(27, 86, 45)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 3).Value = $vals[0]
    $ws.Cells.Item($row, 4).Value = $vals[1]
    $ws.Cells.Item($row, 5).Value = $vals[2]
}

# All three criteria columns now share one integer number format.
$ws.Range("C2:E21").NumberFormat = "0"

# Put the cursor on C2, matching where the user was last working.
[void]$ws.Range("C2").Select()
